$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell H1 "Save" - copy style from G1 (the "sum" header) so it gets
# the same bold/centered/bordered formatting (style index 1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# New Save column values (all 0) for rows 2-7
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
